$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PSSM value updates (supplemental figures recalculation): rows 2-21 x cols B-K
$pssmRows = @(
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, 2.698584711816496, -17.61054187717835, -17.61054187717835),
    @(-17.61054187717835, 4.321921246349989, 1.660303207743798, -17.61054187717835, 3.361961894648361, -17.61054187717835, 1.485764064543614, -17.61054187717835, 2.249697936124011, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, 2.872322197441447, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835),
    @(2.681673018199877, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, 2.104811870583196, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835),
    @(3.763827437036832, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, 1.168064689389561, -17.61054187717835, 1.899126828989616),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, 2.550812702695334, -17.61054187717835, 2.78287057751013, -17.61054187717835, -17.61054187717835, -17.61054187717835, 1.815044898620866),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, 2.450482492389792, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, 1.999208598561589, 1.938928836299899),
    @(-17.61054187717835, -17.61054187717835, 1.391748488582701, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, 2.142983389038862),
    @(-17.61054187717835, -17.61054187717835, 1.504636429889951, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, 2.103168501585743, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, 2.313622569169007, -17.61054187717835, -17.61054187717835, -17.61054187717835, 1.327995363961509, 2.018161742049607, 2.020082535157828, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, -17.61054187717835, 1.595988431156097, 1.262917255321488, 1.532913585354992, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, 1.56703581677363, -17.61054187717835, -17.61054187717835, -17.61054187717835, 1.502755176299086, 1.500820521147286, -17.61054187717835, -17.61054187717835),
    @(-17.61054187717835, -17.61054187717835, 1.783814048287094, -17.61054187717835, 3.280738248270968, -17.61054187717835, 1.987112445291198, 0.9859864077578275, -17.61054187717835, 2.170050752431339),
    @(-17.61054187717835, -17.61054187717835, -17.61054187717835, 2.128581104524647, -17.61054187717835, 2.534879100910471, 2.283169963657213, -17.61054187717835, -17.61054187717835, -17.61054187717835)
)

for ($i = 0; $i -lt $pssmRows.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $pssmRows[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

Write-Output "Updated PSSM values for rows 2-21, columns B-K"
